# Fills in the four empty sub-bullets of the second indicator list
# (Aufgabe 1, "Geben Sie ... an, die Sie als verteilte Anwendung ...")
# with the example text that was added in the commit "so fertig wie geht".
#
# Each target is an empty list paragraph (ilvl=1, numId=16) that directly
# follows a specific labelled bullet ("(A1) ...", "(A3) ...", "(A5) ...",
# "(N3) ..."). Several of those labels also occur earlier in the document
# (in the GIT/version-control table), so we anchor on "label text followed
# by an empty paragraph" rather than on the label text alone.

$d = $word.ActiveDocument

function Find-EmptyParagraphAfter($marker) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text -like "*$marker*") {
            $nextPara = $para.Next()
            if ($nextPara -ne $null) {
                $nextText = $nextPara.Range.Text
                $nextText = $nextText.Trim()
                if ($nextText -eq "") {
                    return $nextPara
                }
            }
        }
    }
    return $null
}

$target_A1 = Find-EmptyParagraphAfter "(A1) gemeinsame Nutzung von Ressourcen"
$target_A1.Range.Text = "Riesenprojekte, wie die Google Suche o.Ä."

$target_A3 = Find-EmptyParagraphAfter "(A3) parallele / nebenläufige Ausführung von Aktivitäten"
$target_A3.Range.Text = "Android Applikationen, die Prozessorkerne für unterschiedliche Dinge verwenden. UI-Thread, DB-Thread etc."

$target_A5 = Find-EmptyParagraphAfter "(A5) erhöhte Fehlertoleranz"
$target_A5.Range.Text = "Services, wie online Videospiele. Wenn ein Server ausfällt, wird ein anderer verwendet."

$target_N3 = Find-EmptyParagraphAfter "(N3) höhere Wartungsaufwände"
$target_N3.Range.Text = "Bei jedem winzigen Projekt, macht es keinen Sinn die Anwendung verteilt aufzubauen, da dies unnötig Arbeit u.A. im Wartungsaufwand mit sich bringt."

Write-Output "A1 -> $($target_A1.Range.Text)"
Write-Output "A3 -> $($target_A3.Range.Text)"
Write-Output "A5 -> $($target_A5.Range.Text)"
Write-Output "N3 -> $($target_N3.Range.Text)"
